# Insert a new weekly price record for Haba (Vega Central Mapocho de Santiago)
# at row 330, pushing the existing rows 330-402 down to 331-403.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 330..402 down by one row.
$ws.Rows(330).Insert()

# Populate the newly inserted row 330 with the new data point.
$ws.Cells.Item(330, 1).Value2  = 9
$ws.Cells.Item(330, 2).Value2  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(330, 3).Value2  = "Metropolitana"
$ws.Cells.Item(330, 4).Value2  = 45211
$ws.Cells.Item(330, 5).Value2  = 13
$ws.Cells.Item(330, 6).Value2  = 100112026
$ws.Cells.Item(330, 7).Value2  = "Haba"
$ws.Cells.Item(330, 8).Value2  = "Sin especificar"
$ws.Cells.Item(330, 9).Value2  = "Primera"
$ws.Cells.Item(330, 10).Value2 = 70
$ws.Cells.Item(330, 11).Value2 = 10000
$ws.Cells.Item(330, 12).Value2 = 11000
$ws.Cells.Item(330, 13).Value2 = 10500
$ws.Cells.Item(330, 14).Value2 = "$/saco 25 kilos"
$ws.Cells.Item(330, 15).Value2 = "Provincia de Melipilla"
$ws.Cells.Item(330, 16).Value2 = 420
$ws.Cells.Item(330, 17).Value2 = 25
$ws.Cells.Item(330, 18).Value2 = "Hortaliza"
